# Update scripts wuth new tpm
# Refresh the NATMI ligand-receptor (Gnas-Tshr) expression/specificity/edge-weight
# values on the active sheet to reflect the recomputed TPM-based inputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 175.411433
$ws.Cells.Item(2, 8).Value = 526.234299
$ws.Cells.Item(2, 9).Value = 0.1535106429347505
$ws.Cells.Item(2, 10).Value = 0.1535106429347505
$ws.Cells.Item(2, 13).Value = 0.6327629999999999
$ws.Cells.Item(2, 14).Value = 1.898289
$ws.Cells.Item(2, 15).Value = 0.1382544270550543
$ws.Cells.Item(2, 16).Value = 0.1382544270550544
$ws.Cells.Item(2, 17).Value = 110.993864579379
$ws.Cells.Item(2, 18).Value = 998.9447812144108
$ws.Cells.Item(2, 19).Value = 0.02122352598579695
$ws.Cells.Item(2, 20).Value = 0.02122352598579696
$ws.Cells.Item(3, 7).Value = 175.411433
$ws.Cells.Item(3, 8).Value = 526.234299
$ws.Cells.Item(3, 9).Value = 0.1535106429347505
$ws.Cells.Item(3, 10).Value = 0.1535106429347505
$ws.Cells.Item(3, 15).Value = 0.4765301499162115
$ws.Cells.Item(3, 16).Value = 0.4765301499162115
$ws.Cells.Item(3, 17).Value = 382.5694703195943
$ws.Cells.Item(3, 18).Value = 3443.125232876349
$ws.Cells.Item(3, 19).Value = 0.07315244969143064
$ws.Cells.Item(3, 20).Value = 0.07315244969143068
$ws.Cells.Item(4, 7).Value = 175.411433
$ws.Cells.Item(4, 8).Value = 526.234299
$ws.Cells.Item(4, 9).Value = 0.1535106429347505
$ws.Cells.Item(4, 10).Value = 0.1535106429347505
$ws.Cells.Item(4, 13).Value = 1.444396333333334
$ws.Cells.Item(4, 14).Value = 4.333189000000001
$ws.Cells.Item(4, 15).Value = 0.3155908096798033
$ws.Cells.Item(4, 16).Value = 0.3155908096798033
$ws.Cells.Item(4, 17).Value = 253.3636306499457
$ws.Cells.Item(4, 18).Value = 2280.272675849511
$ws.Cells.Item(4, 19).Value = 0.04844654809824508
$ws.Cells.Item(4, 20).Value = 0.04844654809824509
$ws.Cells.Item(5, 7).Value = 175.411433
$ws.Cells.Item(5, 8).Value = 526.234299
$ws.Cells.Item(5, 9).Value = 0.1535106429347505
$ws.Cells.Item(5, 10).Value = 0.1535106429347505
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.3186579999999999
$ws.Cells.Item(5, 14).Value = 0.9559739999999999
$ws.Cells.Item(5, 15).Value = 0.06962461334893082
$ws.Cells.Item(5, 16).Value = 0.06962461334893082
$ws.Cells.Item(5, 17).Value = 55.89625641691399
$ws.Cells.Item(5, 18).Value = 503.0663077522259
$ws.Cells.Item(5, 19).Value = 0.01068811915927778
$ws.Cells.Item(5, 20).Value = 0.01068811915927778
$ws.Cells.Item(6, 9).Value = 0.464799214434963
$ws.Cells.Item(6, 10).Value = 0.4647992144349631
$ws.Cells.Item(6, 13).Value = 0.6327629999999999
$ws.Cells.Item(6, 14).Value = 1.898289
$ws.Cells.Item(6, 15).Value = 0.1382544270550543
$ws.Cells.Item(6, 16).Value = 0.1382544270550544
$ws.Cells.Item(6, 17).Value = 336.0669988563869
$ws.Cells.Item(6, 18).Value = 3024.602989707482
$ws.Cells.Item(6, 19).Value = 0.06426054908734516
$ws.Cells.Item(6, 20).Value = 0.06426054908734517
$ws.Cells.Item(7, 9).Value = 0.464799214434963
$ws.Cells.Item(7, 10).Value = 0.4647992144349631
$ws.Cells.Item(7, 15).Value = 0.4765301499162115
$ws.Cells.Item(7, 16).Value = 0.4765301499162115
$ws.Cells.Item(7, 19).Value = 0.2214908393356302
$ws.Cells.Item(7, 20).Value = 0.2214908393356303
$ws.Cells.Item(8, 9).Value = 0.464799214434963
$ws.Cells.Item(8, 10).Value = 0.4647992144349631
$ws.Cells.Item(8, 13).Value = 1.444396333333334
$ws.Cells.Item(8, 14).Value = 4.333189000000001
$ws.Cells.Item(8, 15).Value = 0.3155908096798033
$ws.Cells.Item(8, 16).Value = 0.3155908096798033
$ws.Cells.Item(8, 17).Value = 767.1338888375316
$ws.Cells.Item(8, 18).Value = 6904.204999537784
$ws.Cells.Item(8, 19).Value = 0.1466863604220665
$ws.Cells.Item(8, 20).Value = 0.1466863604220665
$ws.Cells.Item(9, 9).Value = 0.464799214434963
$ws.Cells.Item(9, 10).Value = 0.4647992144349631
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.3186579999999999
$ws.Cells.Item(9, 14).Value = 0.9559739999999999
$ws.Cells.Item(9, 15).Value = 0.06962461334893082
$ws.Cells.Item(9, 16).Value = 0.06962461334893082
$ws.Cells.Item(9, 17).Value = 169.2425722135753
$ws.Cells.Item(9, 18).Value = 1523.183149922178
$ws.Cells.Item(9, 19).Value = 0.03236146558992108
$ws.Cells.Item(9, 20).Value = 0.03236146558992109
$ws.Cells.Item(10, 7).Value = 360.115397
$ws.Cells.Item(10, 8).Value = 1080.346191
$ws.Cells.Item(10, 9).Value = 0.3151536087398187
$ws.Cells.Item(10, 10).Value = 0.3151536087398188
$ws.Cells.Item(10, 13).Value = 0.6327629999999999
$ws.Cells.Item(10, 14).Value = 1.898289
$ws.Cells.Item(10, 15).Value = 0.1382544270550543
$ws.Cells.Item(10, 16).Value = 0.1382544270550544
$ws.Cells.Item(10, 17).Value = 227.867698951911
$ws.Cells.Item(10, 18).Value = 2050.809290567199
$ws.Cells.Item(10, 19).Value = 0.04357138161065641
$ws.Cells.Item(10, 20).Value = 0.04357138161065642
$ws.Cells.Item(11, 7).Value = 360.115397
$ws.Cells.Item(11, 8).Value = 1080.346191
$ws.Cells.Item(11, 9).Value = 0.3151536087398187
$ws.Cells.Item(11, 10).Value = 0.3151536087398188
$ws.Cells.Item(11, 15).Value = 0.4765301499162115
$ws.Cells.Item(11, 16).Value = 0.4765301499162115
$ws.Cells.Item(11, 17).Value = 785.4057989721824
$ws.Cells.Item(11, 18).Value = 7068.652190749642
$ws.Cells.Item(11, 19).Value = 0.1501801964194209
$ws.Cells.Item(11, 20).Value = 0.1501801964194209
$ws.Cells.Item(12, 7).Value = 360.115397
$ws.Cells.Item(12, 8).Value = 1080.346191
$ws.Cells.Item(12, 9).Value = 0.3151536087398187
$ws.Cells.Item(12, 10).Value = 0.3151536087398188
$ws.Cells.Item(12, 13).Value = 1.444396333333334
$ws.Cells.Item(12, 14).Value = 4.333189000000001
$ws.Cells.Item(12, 15).Value = 0.3155908096798033
$ws.Cells.Item(12, 16).Value = 0.3155908096798033
$ws.Cells.Item(12, 17).Value = 520.1493590036778
$ws.Cells.Item(12, 18).Value = 4681.3442310331
$ws.Cells.Item(12, 19).Value = 0.09945958255571133
$ws.Cells.Item(12, 20).Value = 0.09945958255571134
$ws.Cells.Item(13, 7).Value = 360.115397
$ws.Cells.Item(13, 8).Value = 1080.346191
$ws.Cells.Item(13, 9).Value = 0.3151536087398187
$ws.Cells.Item(13, 10).Value = 0.3151536087398188
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.3186579999999999
$ws.Cells.Item(13, 14).Value = 0.9559739999999999
$ws.Cells.Item(13, 15).Value = 0.06962461334893082
$ws.Cells.Item(13, 16).Value = 0.06962461334893082
$ws.Cells.Item(13, 17).Value = 114.753652177226
$ws.Cells.Item(13, 18).Value = 1032.782869595034
$ws.Cells.Item(13, 19).Value = 0.0219424481540301
$ws.Cells.Item(13, 20).Value = 0.02194244815403011
$ws.Cells.Item(14, 7).Value = 76.02905266666666
$ws.Cells.Item(14, 8).Value = 228.087158
$ws.Cells.Item(14, 9).Value = 0.06653653389046771
$ws.Cells.Item(14, 10).Value = 0.06653653389046772
$ws.Cells.Item(14, 13).Value = 0.6327629999999999
$ws.Cells.Item(14, 14).Value = 1.898289
$ws.Cells.Item(14, 15).Value = 0.1382544270550543
$ws.Cells.Item(14, 16).Value = 0.1382544270550544
$ws.Cells.Item(14, 17).Value = 48.10837145251799
$ws.Cells.Item(14, 18).Value = 432.9753430726619
$ws.Cells.Item(14, 19).Value = 0.009198970371255819
$ws.Cells.Item(14, 20).Value = 0.009198970371255822
$ws.Cells.Item(15, 7).Value = 76.02905266666666
$ws.Cells.Item(15, 8).Value = 228.087158
$ws.Cells.Item(15, 9).Value = 0.06653653389046771
$ws.Cells.Item(15, 10).Value = 0.06653653389046772
$ws.Cells.Item(15, 15).Value = 0.4765301499162115
$ws.Cells.Item(15, 16).Value = 0.4765301499162115
$ws.Cells.Item(15, 17).Value = 165.8181220581397
$ws.Cells.Item(15, 18).Value = 1492.363098523258
$ws.Cells.Item(15, 19).Value = 0.03170666446972966
$ws.Cells.Item(15, 20).Value = 0.03170666446972967
$ws.Cells.Item(16, 7).Value = 76.02905266666666
$ws.Cells.Item(16, 8).Value = 228.087158
$ws.Cells.Item(16, 9).Value = 0.06653653389046771
$ws.Cells.Item(16, 10).Value = 0.06653653389046772
$ws.Cells.Item(16, 13).Value = 1.444396333333334
$ws.Cells.Item(16, 14).Value = 4.333189000000001
$ws.Cells.Item(16, 15).Value = 0.3155908096798033
$ws.Cells.Item(16, 16).Value = 0.3155908096798033
$ws.Cells.Item(16, 17).Value = 109.8160848985402
$ws.Cells.Item(16, 18).Value = 988.3447640868621
$ws.Cells.Item(16, 19).Value = 0.02099831860378037
$ws.Cells.Item(16, 20).Value = 0.02099831860378038
$ws.Cells.Item(17, 7).Value = 76.02905266666666
$ws.Cells.Item(17, 8).Value = 228.087158
$ws.Cells.Item(17, 9).Value = 0.06653653389046771
$ws.Cells.Item(17, 10).Value = 0.06653653389046772
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.3186579999999999
$ws.Cells.Item(17, 14).Value = 0.9559739999999999
$ws.Cells.Item(17, 15).Value = 0.06962461334893082
$ws.Cells.Item(17, 16).Value = 0.06962461334893082
$ws.Cells.Item(17, 17).Value = 24.22726586465466
$ws.Cells.Item(17, 18).Value = 218.045392781892
$ws.Cells.Item(17, 19).Value = 0.004632580445701846
$ws.Cells.Item(17, 20).Value = 0.004632580445701847
